# Regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Only the "K" column (column G) values change for rows 2-23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 0
    6  = 0
    7  = 0
    8  = 1
    9  = 0
    10 = 1
    11 = 0
    12 = 1
    13 = 0
    14 = 0
    15 = 2
    16 = 1
    17 = 0
    18 = 0
    19 = 1
    20 = 3
    21 = 1
    22 = 2
    23 = 0
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
